# Generate Report for Archive
# Update the handoff status text everywhere it appears ("Ready for handoff"
# -> "In Translation") and let the affected Status columns re-fit to the
# new (shorter) text width.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ("Ready for handoff" -eq $cell.Value2) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# Re-fit the Status column(s) now that the text is shorter.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
